$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3 = 0
    4 = 0
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 0
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 0
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 2
    20 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 2
    27 = 4
    28 = 2
    29 = 1
    30 = 1
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 2
    38 = 0
    39 = 1
    40 = 0
    41 = 1
    42 = 1
    43 = 0
    44 = 0
    45 = 1
    46 = 3
    47 = 1
    48 = 0
    49 = 4
    50 = 2
    51 = 3
    52 = 0
    53 = 0
    54 = 0
    55 = 2
    56 = 1
    57 = 2
    58 = 2
    59 = 0
    60 = 1
    61 = 2
    62 = 1
    63 = 1
    64 = 2
    65 = 2
    66 = 1
    67 = 0
    68 = 1
    69 = 2
    71 = 1
    72 = 0
    73 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
